$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '42.645.41'
Set-TextValue 'E2' '  +2.91%  '
Set-TextValue 'D3' '2.207.83'
Set-TextValue 'E3' '  +2.27%  '
Set-TextValue 'E4' '  -0.18%  '
Set-TextValue 'D5' '249.79'
Set-TextValue 'E5' '  +5.38%  '
Set-TextValue 'D6' '0.613'
Set-TextValue 'E6' '  +1.46%  '
Set-TextValue 'D7' '74.37'
Set-TextValue 'E7' '  +5.08%  '
Set-TextValue 'E8' '  -0.12%  '
Set-TextValue 'D9' '0.589'
Set-TextValue 'E9' '  +2.74%  '
Set-TextValue 'D10' '40.45'
Set-TextValue 'E10' '  +1.80%  '
Set-TextValue 'D11' '0.0914'
Set-TextValue 'E11' '  +1.13%  '
Set-TextValue 'D12' '6.82'
Set-TextValue 'E12' '  +2.57%  '
Set-TextValue 'D13' '0.101'
Set-TextValue 'E13' '  +1.26%  '
Set-TextValue 'D14' '2.542.69'
Set-TextValue 'E14' '  +2.46%  '
Set-TextValue 'D15' '14.39'
Set-TextValue 'E15' '  +0.63%  '
Set-TextValue 'D16' '2.213.76'
Set-TextValue 'E16' '  +3.34%  '
Set-TextValue 'D17' '0.779'
Set-TextValue 'E17' '  -0.46%  '
Set-TextValue 'D18' '42.553.47'
Set-TextValue 'E18' '  +3.11%  '
Set-TextValue 'D19' '0.0000102'
Set-TextValue 'E19' '  +1.76%  '
Set-TextValue 'D20' '70.97'
Set-TextValue 'E20' '  +2.38%  '
Set-TextValue 'D21' '5.89'
Set-TextValue 'E21' '  +2.53%  '
Set-TextValue 'D22' '228.74'
Set-TextValue 'E22' '  +1.57%  '
Set-TextValue 'D23' '2.16'
Set-TextValue 'E23' '  +9.84%  '
Set-TextValue 'D24' '9.40'
Set-TextValue 'E24' '  -3.69%  '
Set-TextValue 'E25' '  -0.01%  '
Set-TextValue 'D26' '10.69'
Set-TextValue 'E26' '  +0.58%  '
Set-TextValue 'D27' '3.40'
Set-TextValue 'E27' '  +2.13%  '
Set-TextValue 'D28' '38.50'
Set-TextValue 'E28' '  +21.47%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D29' '2.23'
Set-TextValue 'E29' '  +2.71%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D30' '2.20'
Set-TextValue 'E30' '  +1.48%  '
Set-TextValue 'D31' '169.60'
Set-TextValue 'E31' '  -0.98%  '
Set-TextValue 'D32' '20.08'
Set-TextValue 'E32' '  +1.83%  '
Set-TextValue 'D33' '0.0788'
Set-TextValue 'E33' '  +3.17%  '
Set-TextValue 'D34' '5.18'
Set-TextValue 'E34' '  +1.85%  '
Set-TextValue 'D35' '0.120'
Set-TextValue 'E35' '  +0.05%  '
Set-TextValue 'D36' '0.109'
Set-TextValue 'E36' '  +4.73%  '
Set-TextValue 'D37' '4.39'
Set-TextValue 'E37' '  +1.66%  '
Set-TextValue 'E38' '  +10.60%  '
Set-TextValue 'D39' '11.99'
Set-TextValue 'E39' '  -1.31%  '
Set-TextValue 'D40' '2.08'
Set-TextValue 'E40' '  +1.35%  '
Set-TextValue 'D41' '0.201'
Set-TextValue 'E41' '  +7.00%  '
Set-TextValue 'D42' '5.27'
Set-TextValue 'E42' '  -0.85%  '
Set-TextValue 'D43' '58.64'
Set-TextValue 'E43' '  +0.74%  '
Set-TextValue 'D44' '8.50'
Set-TextValue 'E44' '  +3.07%  '
Set-TextValue 'D45' '102.11'
Set-TextValue 'E45' '  +4.62%  '
$ws.Range('B46').Value = 'WOONetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue 'D46' '0.472'
Set-TextValue 'E46' '  +21.22%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D47' '0.0978'
Set-TextValue 'E47' '  +1.81%  '
Set-TextValue 'D48' '2.40'
Set-TextValue 'E48' '  +11.17%  '
Set-TextValue 'D49' '1.10'
Set-TextValue 'E49' '  +2.55%  '
$ws.Range('B50').Value = 'FTXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D50' '4.15'
Set-TextValue 'E50' '  +20.20%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D51' '1.12'
Set-TextValue 'E51' '  +1.32%  '

Write-Host "Applied cryptos update."
